$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename distribution type labels to match the values found in the enum
$ws.Range("F4").Value = "PERT"
$ws.Range("F5").Value = "LogUniform"
$ws.Range("F6").Value = "TruncatedNormal"
$ws.Range("F7").Value = "TruncatedLogNormal"

# Update the active selection on the sheet
$ws.Range("G10").Select()
